$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before E (shifts current "note" column E -> F)
$ws.Range("E1").EntireColumn.Insert()

# New header cell E4 "step "
$ws.Range("E4").Value = "step "

# New numeric value E7
$ws.Range("E7").Value = 75

# New row 8 data
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A8").Value = 45661
$ws.Range("B8").Value = -0.14788000000000001
$ws.Range("D8").Value = -0.30982999999999999
$ws.Range("E8").Value = 75
$ws.Range("F8").Value = "pid tuner -- peak: -0.6 theta, 20 omega, 0.3 torque"

# Update selection to mirror final state
$ws.Range("I11").Select()
